$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Outputs")
$ws.Activate()

# --- Renumber the "index" column (B2:B17) down by one (0-based indexing) ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 2
}

# --- Move / resize the embedded chart so it spans roughly E22:O53 instead of B20:E50 ---
$co = $ws.ChartObjects(1)
$co.Left   = 769.8581160802165
$co.Top    = 313.2605511811024
$co.Width  = 607.4541264763781
$co.Height = 437.85456692913385

# --- Update the view: zoom to 99% and change the active selection ---
$excel.ActiveWindow.Zoom = 99
[void]$ws.Range("B14:B17").Select()
